$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column B (TOTAL_SUBSTATION_LOAD), C (CONTESTABLE_ENERGY), D (ACTUAL_ENERGY) updates per row.
# Rows 2-19 get new B, C, D values.
# Rows 20-25 have their B cell cleared (no TOTAL_SUBSTATION_LOAD data anymore) and get new C values;
# D values also change for rows 20-24 (row 25's D stays the same).

$data = @(
    @{ Row = 2;  B = 72056;  C = 5516.6335;            D = 66539.3665 }
    @{ Row = 3;  B = 69308;  C = 5387.312;              D = 63920.688 }
    @{ Row = 4;  B = 65612;  C = 5339.824000000001;     D = 60272.176 }
    @{ Row = 5;  B = 64269;  C = 5306.819000000001;     D = 58962.181 }
    @{ Row = 6;  B = 65136;  C = 5338.83;                D = 59797.17 }
    @{ Row = 7;  B = 71065;  C = 5422.714499999999;     D = 65642.2855 }
    @{ Row = 8;  B = 70049;  C = 5875.107;               D = 64173.893 }
    @{ Row = 9;  B = 80484;  C = 7214.4275;              D = 73269.57249999999 }
    @{ Row = 10; B = 94725;  C = 8467.018;               D = 86257.982 }
    @{ Row = 11; B = 104123; C = 13640.6165;             D = 90482.3835 }
    @{ Row = 12; B = 110312; C = 15763.02;               D = 94548.98 }
    @{ Row = 13; B = 111321; C = 15072.0395;             D = 96248.9605 }
    @{ Row = 14; B = 109259; C = 15259.6675;             D = 93999.3325 }
    @{ Row = 15; B = 113429; C = 15815.744;              D = 97613.25599999999 }
    @{ Row = 16; B = 114542; C = 15955.3275;             D = 98586.6725 }
    @{ Row = 17; B = 99426;  C = 16114.0735;             D = 83311.9265 }
    @{ Row = 18; B = 92477;  C = 16384.389;              D = 76092.611 }
    @{ Row = 19; B = 1237;   C = 15932.637;              D = 0 }
)

foreach ($item in $data) {
    $r = $item.Row
    $ws.Cells.Item($r, 2).Value = $item.B
    $ws.Cells.Item($r, 3).Value = $item.C
    $ws.Cells.Item($r, 4).Value = $item.D
}

# Rows 20-25: clear column B, update column C, and update column D (row 25's D is unchanged).
$dataNoB = @(
    @{ Row = 20; C = 15036.9065;             D = 81334.0935 }
    @{ Row = 21; C = 13330.5025;             D = 78578.4975 }
    @{ Row = 22; C = 11770.8115;             D = 77749.1885 }
    @{ Row = 23; C = 9312.075499999999;      D = 74837.92449999999 }
    @{ Row = 24; C = 6773.7985;              D = 54075.2015 }
    @{ Row = 25; C = 5542.130999999999;      D = 14561.85 }
)

foreach ($item in $dataNoB) {
    $r = $item.Row
    $ws.Cells.Item($r, 2).ClearContents()
    $ws.Cells.Item($r, 3).Value = $item.C
    $ws.Cells.Item($r, 4).Value = $item.D
}
